{"js": "// Commit: \"Minor change to list - removed 'and'\"\n//\n// In the \"Skills and Interests\" -> \"Technical Skills\" table, the Programming\n// and Scripting Languages list drops the word \"and\" before \"Bash\", and the\n// Tools list drops the word \"and\" before \"TensorFlow.\" (both were Oxford-\n// comma \"and\"s joining the final item of a comma list).\n//\n//   \"R, Python, Java, JavaScript, Go, and Bash\"        -> \"R, Python, Java, JavaScript, Go, Bash\"\n//   \"Cloud Platform, Theano, and TensorFlow.\"           -> \"Cloud Platform, Theano, TensorFlow.\"\n\nconst body = context.document.body;\n\n// 1) \"...Go, and Bash\" -> \"...Go, Bash\"\nconst langResults = body.search(\"R, Python, Java, JavaScript, Go, and Bash\", { matchCase: true });\nlangResults.load(\"items\");\nawait context.sync();\n\nif (langResults.items.length > 0) {\n  langResults.items[0].insertText(\n    \"R, Python, Java, JavaScript, Go, Bash\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2) \"Theano, and TensorFlow\" -> \"Theano, TensorFlow\"\nconst toolsResults = body.search(\"Cloud Platform, Theano, and TensorFlow\", { matchCase: true });\ntoolsResults.load(\"items\");\nawait context.sync();\n\nif (toolsResults.items.length > 0) {\n  toolsResults.items[0].insertText(\n    \"Cloud Platform, Theano, TensorFlow\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Commit: \"Minor change to list - removed 'and'\"\n#\n# In the \"Skills and Interests\" -> \"Technical Skills\" table, the Programming\n# and Scripting Languages list drops the word \"and\" before \"Bash\", and the\n# Tools list drops the word \"and\" before \"TensorFlow.\" (both were Oxford-\n# comma \"and\"s joining the final item of a comma list).\n#\n#   \"R, Python, Java, JavaScript, Go, and Bash\"  -> \"R, Python, Java, JavaScript, Go, Bash\"\n#   \"Cloud Platform, Theano, and TensorFlow.\"     -> \"Cloud Platform, Theano, TensorFlow.\"\n\n$d = $word.ActiveDocument\n\n# 1) \"...Go, and Bash\" -> \"...Go, Bash\"\n$r1 = $d.Content\n$r1.Find.Execute(\"R, Python, Java, JavaScript, Go, and Bash\")\nif ($r1.Find.Found) {\n    $r1.Text = \"R, Python, Java, JavaScript, Go, Bash\"\n}\n\n# 2) \"Theano, and TensorFlow\" -> \"Theano, TensorFlow\"\n$r2 = $d.Content\n$r2.Find.Execute(\"Cloud Platform, Theano, and TensorFlow\")\nif ($r2.Find.Found) {\n    $r2.Text = \"Cloud Platform, Theano, TensorFlow\"\n}\n"}
